$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("XPath")

# 1. Update the existing "payment details tab" xpath (B63):
#    '(//LI)[23]' -> "(//SPAN[text()='Payment Details'])[1]"
$ws.Range("B63").Value = "(//SPAN[text()='Payment Details'])[1]"

# 2. New row 64: current_address / DIV xpath for the current address
$ws.Rows.Item(64).RowHeight = 28.2
$ws.Range("A64").Value = "current_address"
$ws.Range("B64").Value = "(//DIV[text()='Electronic City phase 1 , Konappana agrahara ,nanjuda reddy layout near yellamma temple ,Hosur main road Bangalore'])[1]"

# B64 should look like the other "xpath value" cells (bold green Menlo) - copy that
# format from an existing xpath cell (B60) onto both B64 and A64.
$ws.Range("B60").Copy()
$ws.Range("B64").PasteSpecial(-4122)
$ws.Range("A64").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# A64 (the label column) keeps the Menlo/12pt look but in plain black, not bold -
# tweak the pasted font accordingly. This produces the new font + cell style used
# by the other new label cells below (A65, A66).
$ws.Range("A64").Font.Color = 0
$ws.Range("A64").Font.FontStyle = "Regular"

# 3. New row 65: address_change_dt / DTButton xpath [1]
$ws.Rows.Item(65).RowHeight = 28.2
$ws.Range("A65").Value = "address_change_dt"
$ws.Range("B65").Value = "(//DIV[@class='_-sa-flipkart-src-Components-DT-DTButton-DTButton_button _-sa-flipkart-src-Components-DT-DTButton-DTButton_outlineBtn _-sa-flipkart-src-Components-DT-DTButton-DTButton_isDisabled'])[1]"

# 4. New row 66: create_incident_dt / DTButton xpath [2]
$ws.Rows.Item(66).RowHeight = 28.2
$ws.Range("A66").Value = "create_incident_dt"
$ws.Range("B66").Value = "(//DIV[@class='_-sa-flipkart-src-Components-DT-DTButton-DTButton_button _-sa-flipkart-src-Components-DT-DTButton-DTButton_outlineBtn _-sa-flipkart-src-Components-DT-DTButton-DTButton_isDisabled'])[2]"

# Give A65 / A66 the same label font as A64
$ws.Range("A64").Copy()
$ws.Range("A65").PasteSpecial(-4122)
$ws.Range("A66").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 5. New row 67: address tab / li[3] image xpath (plain default style, default row height)
$ws.Range("A67").Value = "address tab"
$ws.Range("B67").Value = '//*[@id="root"]/div/div[4]/div[1]/div/div/div/div[2]/div[1]/div[3]/div[2]/div/div[1]/ul/li[3]/div/div/img'

# 6. Scroll / selection bookkeeping to match the edited view
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 55
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B71").Select()

Write-Output "edit applied"
